# Add three new auction-listing rows (6, 7, 8) to Sheet1, matching the
# newly scraped "Outros" properties from leilaoimovel.com.br / Caninde-CE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Assign the value, then reset the style to "Normal" so Excel does not
    # keep an implicit Text/quote-prefix number format on the cell (which
    # would otherwise add an s="..." style attribute not present in the
    # original template rows).
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).Style = "Normal"
}

# ---------------------------------------------------------------- Row 6 --
Set-TextCell "A6" " Outros "
Set-TextCell "K6" "'05/03/2024"
Set-TextCell "O6" "Judicial"
Set-TextCell "P6" " Norte Nordeste Leilões "
$ws.Range("R6").Value = 0
Set-TextCell "S6" "R. RAIMUNDO ALCONFORADO, Nº 158"
Set-TextCell "T6" "https://www.leilaoimovel.com.br/imovel/ce/caninde/outros-predio-2-andares-terreno-200m-caninde-ce-desocupado-imovel-1627069"

# ---------------------------------------------------------------- Row 7 --
Set-TextCell "A7" " Outros "
Set-TextCell "K7" "'05/03/2024"
Set-TextCell "O7" "Judicial"
Set-TextCell "P7" " Norte Nordeste Leilões "
$ws.Range("R7").Value = 0
Set-TextCell "S7" "RUA JOSÉ VELOSO JUCÁ, N° 2576"
Set-TextCell "T7" "https://www.leilaoimovel.com.br/imovel/ce/caninde/outros-predio-triplex-terreno-130m-caninde-ce-imovel-1627070"

# ---------------------------------------------------------------- Row 8 --
Set-TextCell "A8" " Outros "
Set-TextCell "K8" "'05/03/2024"
Set-TextCell "O8" "Judicial"
Set-TextCell "P8" " Norte Nordeste Leilões "
$ws.Range("R8").Value = 0
Set-TextCell "S8" "TRAVESSA JOÃO MARTINS, N° 56"
Set-TextCell "T8" "https://www.leilaoimovel.com.br/imovel/ce/caninde/outros-predio-2-pavs-196-20m-terreno-150m-caninde-ce-imovel-1627068"

Write-Host "Added rows 6-8. New dimension should be A1:T8."
